$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1219.3889
$ws.Range("J17").Value = 1238.7646
$ws.Range("L17").Value = 3716.2938
$ws.Range("N17").Value = -4052.2938

$ws.Range("H64").Value = 3424.6365
$ws.Range("I64").Value = 3061.0715
$ws.Range("J64").Value = 3692.5264
$ws.Range("K64").Value = 3061.0715
$ws.Range("L64").Value = 3692.5264
$ws.Range("M64").Value = -2813.0715
$ws.Range("N64").Value = -4188.526400000001

$ws.Range("H67").Value = 3424.6365
$ws.Range("I67").Value = 3061.0715
$ws.Range("J67").Value = 3692.5264
$ws.Range("K67").Value = 3061.0715
$ws.Range("L67").Value = 3692.5264
$ws.Range("M67").Value = -2203.0715
$ws.Range("N67").Value = -5408.526400000001

$ws.Range("H76").Value = 3177.6
$ws.Range("I76").Value = 2828
$ws.Range("J76").Value = 3993.3333
$ws.Range("K76").Value = 2828
$ws.Range("L76").Value = 3993.3333
$ws.Range("M76").Value = -2513
$ws.Range("N76").Value = -4623.3333

$ws.Range("H79").Value = 3177.6
$ws.Range("I79").Value = 2828
$ws.Range("J79").Value = 3993.3333
$ws.Range("K79").Value = 2828
$ws.Range("L79").Value = 3993.3333
$ws.Range("M79").Value = -1736
$ws.Range("N79").Value = -6177.3333

$ws.Range("H98").Value = 2217.6785
$ws.Range("I98").Value = 2220.4167
$ws.Range("J98").Value = 2201.25
$ws.Range("K98").Value = 2220.4167
$ws.Range("L98").Value = 2201.25
$ws.Range("M98").Value = -722.4167000000002
$ws.Range("N98").Value = -5197.25

$ws.Range("H112").Value = 3029.0322
$ws.Range("I112").Value = 1466.6666
$ws.Range("K112").Value = 4399.9998
$ws.Range("M112").Value = -3291.9998

$ws.Range("H122").Value = 2217.6785
$ws.Range("I122").Value = 2220.4167
$ws.Range("J122").Value = 2201.25
$ws.Range("K122").Value = 6661.250100000001
$ws.Range("L122").Value = 6603.75
$ws.Range("M122").Value = -4211.250100000001
$ws.Range("N122").Value = -11503.75

$ws.Range("H129").Value = 1786.6471
$ws.Range("I129").Value = 631.55554
$ws.Range("J129").Value = 2202.48
$ws.Range("K129").Value = 1894.66662
$ws.Range("L129").Value = 6607.440000000001
$ws.Range("M129").Value = 3105.33338
$ws.Range("N129").Value = -16607.44

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1639377.5
$ws.Range("I32").Value = 13387.623
$ws.Range("J32").Value = 13021307
$ws.Range("K32").Value = 13387.623
$ws.Range("L32").Value = 13021307
$ws.Range("M32").Value = -13100.623
$ws.Range("N32").Value = -13021881

$ws.Range("H61").Value = 2810.7144
$ws.Range("I61").Value = 3304.5
$ws.Range("J61").Value = 2152.3333
$ws.Range("K61").Value = 3304.5
$ws.Range("L61").Value = 2152.3333
$ws.Range("M61").Value = -3092.5
$ws.Range("N61").Value = -2576.3333

$ws.Range("H122").Value = 1227.3334
$ws.Range("I122").Value = 1052.1666
$ws.Range("J122").Value = 1314.9166
$ws.Range("K122").Value = 3156.4998
$ws.Range("L122").Value = 3944.7498
$ws.Range("M122").Value = -706.4998000000001
$ws.Range("N122").Value = -8844.7498

$ws.Range("H136").Value = 2810.7144
$ws.Range("I136").Value = 3304.5
$ws.Range("J136").Value = 2152.3333
$ws.Range("K136").Value = 9913.5
$ws.Range("L136").Value = 6456.999899999999
$ws.Range("M136").Value = -7363.5
$ws.Range("N136").Value = -11556.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9593.4
$ws.Range("J4").Value = 9593.4
$ws.Range("L4").Value = 9593.4
$ws.Range("N4").Value = -9817.4

$ws.Range("H25").Value = 25679
$ws.Range("I25").Value = 7011
$ws.Range("K25").Value = 7011
$ws.Range("M25").Value = -6837

$ws.Range("H31").Value = 2785.9692
$ws.Range("I31").Value = 1697.2391
$ws.Range("J31").Value = 5421.8423
$ws.Range("K31").Value = 1697.2391
$ws.Range("L31").Value = 5421.8423
$ws.Range("M31").Value = -1402.2391
$ws.Range("N31").Value = -6011.8423

$ws.Range("H34").Value = 2785.9692
$ws.Range("I34").Value = 1697.2391
$ws.Range("J34").Value = 5421.8423
$ws.Range("K34").Value = 1697.2391
$ws.Range("L34").Value = 5421.8423
$ws.Range("M34").Value = -1495.2391
$ws.Range("N34").Value = -5825.8423

$ws.Range("H53").Value = 30000
$ws.Range("J53").Value = 30000
$ws.Range("L53").Value = 30000
$ws.Range("N53").Value = -31214

$ws.Range("H96").Value = 17813.545
$ws.Range("J96").Value = 17813.545
$ws.Range("L96").Value = 17813.545
$ws.Range("N96").Value = -23305.545

$ws.Range("H105").Value = 744.6667
$ws.Range("I105").Value = 726.4286
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 726.4286
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 1020.5714
$ws.Range("N105").Value = -4494

$ws.Range("H122").Value = 58824764
$ws.Range("I122").Value = 71429450
$ws.Range("J122").Value = 2901.3333
$ws.Range("K122").Value = 214288350
$ws.Range("L122").Value = 8703.999899999999
$ws.Range("M122").Value = -214285900
$ws.Range("N122").Value = -13603.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3807.762
$ws.Range("I134").Value = 1715
$ws.Range("K134").Value = 5145
$ws.Range("M134").Value = -75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 10923.2
$ws.Range("I122").Value = 19179
$ws.Range("J122").Value = 3699.375
$ws.Range("K122").Value = 57537
$ws.Range("L122").Value = 11098.125
$ws.Range("M122").Value = -55087
$ws.Range("N122").Value = -15998.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 32261444
$ws.Range("I40").Value = 41669748
$ws.Range("J40").Value = 4397.857
$ws.Range("K40").Value = 41669748
$ws.Range("L40").Value = 4397.857
$ws.Range("M40").Value = -41669612
$ws.Range("N40").Value = -4669.857

$ws.Range("H46").Value = 833910.8
$ws.Range("J46").Value = 2500732.5
$ws.Range("L46").Value = 2500732.5
$ws.Range("N46").Value = -2501108.5

$ws.Range("H122").Value = 2717.8718
$ws.Range("I122").Value = 2519.4092
$ws.Range("J122").Value = 2974.7058
$ws.Range("K122").Value = 7558.2276
$ws.Range("L122").Value = 8924.117400000001
$ws.Range("M122").Value = -5108.2276
$ws.Range("N122").Value = -13824.1174

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 173390000
$ws.Range("J2").Value = 200068000
$ws.Range("L2").Value = 200068000
$ws.Range("N2").Value = -200068224

$ws.Range("H122").Value = 2071.7058
$ws.Range("I122").Value = 2384.9167
$ws.Range("J122").Value = 1320
$ws.Range("K122").Value = 7154.750100000001
$ws.Range("L122").Value = 3960
$ws.Range("M122").Value = -4704.750100000001
$ws.Range("N122").Value = -8860
